$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings with rich-text runs) ---
$ws.Cells.Item(8,1).Value = "Volume 32   Number  38"
$ws.Cells.Item(9,3).Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Numeric cell updates (Crime Complaints table) ---
$ws.Cells.Item(15,4).Value = 1
$ws.Cells.Item(15,7).Value = 3
$ws.Cells.Item(15,8).Value = -100
$ws.Cells.Item(15,10).Value = 7
$ws.Cells.Item(15,11).Value = 57.142857142857
$ws.Cells.Item(15,12).Value = 57.142857142857
$ws.Cells.Item(15,13).Value = 83.333333333333
$ws.Cells.Item(15,14).Value = -52.173913043478
$ws.Cells.Item(16,4).Value = 3
$ws.Cells.Item(16,5).Value = -100
$ws.Cells.Item(16,6).Value = 5
$ws.Cells.Item(16,7).Value = 7
$ws.Cells.Item(16,8).Value = -28.571428571428
$ws.Cells.Item(16,10).Value = 50
$ws.Cells.Item(16,11).Value = -8
$ws.Cells.Item(16,12).Value = -6.122448979591
$ws.Cells.Item(16,13).Value = -38.666666666666
$ws.Cells.Item(16,14).Value = -88.557213930348
$ws.Cells.Item(17,4).Value = 2
$ws.Cells.Item(17,5).Value = -100
$ws.Cells.Item(17,6).Value = 4
$ws.Cells.Item(17,7).Value = 6
$ws.Cells.Item(17,8).Value = -33.333333333333
$ws.Cells.Item(17,10).Value = 88
$ws.Cells.Item(17,11).Value = 19.318181818181
$ws.Cells.Item(17,12).Value = 0
$ws.Cells.Item(17,13).Value = 40
$ws.Cells.Item(17,14).Value = -49.760765550239
$ws.Cells.Item(18,4).Value = 2
$ws.Cells.Item(18,5).Value = -100
$ws.Cells.Item(18,6).Value = 4
$ws.Cells.Item(18,7).Value = 14
$ws.Cells.Item(18,8).Value = -71.428571428571
$ws.Cells.Item(18,10).Value = 69
$ws.Cells.Item(18,11).Value = -21.739130434782
$ws.Cells.Item(18,12).Value = -27.027027027027
$ws.Cells.Item(18,13).Value = -70.652173913043
$ws.Cells.Item(18,14).Value = -93.25
$ws.Cells.Item(19,3).Value = 8
$ws.Cells.Item(19,5).Value = 0
$ws.Cells.Item(19,7).Value = 25
$ws.Cells.Item(19,8).Value = 28
$ws.Cells.Item(19,9).Value = 272
$ws.Cells.Item(19,10).Value = 275
$ws.Cells.Item(19,11).Value = -1.090909090909
$ws.Cells.Item(19,12).Value = -14.733542319749
$ws.Cells.Item(19,13).Value = 19.823788546255
$ws.Cells.Item(19,14).Value = -15.264797507788
$ws.Cells.Item(20,3).Value = 4
$ws.Cells.Item(20,5).Value = 100
$ws.Cells.Item(20,6).Value = 9
$ws.Cells.Item(20,7).Value = 10
$ws.Cells.Item(20,8).Value = -10
$ws.Cells.Item(20,9).Value = 96
$ws.Cells.Item(20,10).Value = 146
$ws.Cells.Item(20,11).Value = -34.246575342465
$ws.Cells.Item(20,12).Value = 9.090909090909
$ws.Cells.Item(20,13).Value = -14.285714285714
$ws.Cells.Item(20,14).Value = -93.108399138549
$ws.Cells.Item(21,3).Value = 12
$ws.Cells.Item(21,4).Value = 18
$ws.Cells.Item(21,5).Value = -33.333333333333
$ws.Cells.Item(21,6).Value = 54
$ws.Cells.Item(21,7).Value = 65
$ws.Cells.Item(21,8).Value = -16.923076923076
$ws.Cells.Item(21,9).Value = 584
$ws.Cells.Item(21,10).Value = 637
$ws.Cells.Item(21,11).Value = -8.320251177394
$ws.Cells.Item(21,12).Value = -9.316770186335
$ws.Cells.Item(21,13).Value = -14.117647058823
$ws.Cells.Item(21,14).Value = -81.472081218274
$ws.Cells.Item(24,3).Value = 15
$ws.Cells.Item(24,4).Value = 28
$ws.Cells.Item(24,5).Value = -46.428571428571
$ws.Cells.Item(24,6).Value = 80
$ws.Cells.Item(24,7).Value = 111
$ws.Cells.Item(24,8).Value = -27.927927927927
$ws.Cells.Item(24,9).Value = 739
$ws.Cells.Item(24,10).Value = 1055
$ws.Cells.Item(24,11).Value = -29.952606635071
$ws.Cells.Item(24,12).Value = -33.061594202898
$ws.Cells.Item(24,13).Value = -12.440758293838
$ws.Cells.Item(25,3).Value = 4
$ws.Cells.Item(25,4).Value = 16
$ws.Cells.Item(25,5).Value = -75
$ws.Cells.Item(25,7).Value = 61
$ws.Cells.Item(25,8).Value = -42.622950819672
$ws.Cells.Item(25,9).Value = 341
$ws.Cells.Item(25,10).Value = 659
$ws.Cells.Item(25,11).Value = -48.254931714719
$ws.Cells.Item(25,12).Value = -40.280210157618
$ws.Cells.Item(26,3).Value = 6
$ws.Cells.Item(26,4).Value = 5
$ws.Cells.Item(26,5).Value = 20
$ws.Cells.Item(26,6).Value = 27
$ws.Cells.Item(26,7).Value = 26
$ws.Cells.Item(26,8).Value = 3.846153846153
$ws.Cells.Item(26,9).Value = 293
$ws.Cells.Item(26,10).Value = 261
$ws.Cells.Item(26,11).Value = 12.260536398467
$ws.Cells.Item(26,12).Value = 8.518518518518
$ws.Cells.Item(26,13).Value = 6.934306569343
$ws.Cells.Item(27,4).Value = 1
$ws.Cells.Item(27,7).Value = 4
$ws.Cells.Item(27,8).Value = -100
$ws.Cells.Item(27,10).Value = 10
$ws.Cells.Item(27,11).Value = 50
$ws.Cells.Item(27,12).Value = 50
$ws.Cells.Item(28,6).Value = 4
$ws.Cells.Item(28,7).Value = 5
$ws.Cells.Item(28,8).Value = -20
$ws.Cells.Item(28,9).Value = 31
$ws.Cells.Item(28,11).Value = -26.190476190476
$ws.Cells.Item(28,12).Value = 34.782608695652
$ws.Cells.Item(31,12).Value = -20
$ws.Cells.Item(33,4).Value = 1
$ws.Cells.Item(33,10).Value = 9
$ws.Cells.Item(33,11).Value = -88.888888888888

# --- Cells converted from numeric to text ("0" / "***.*") shared strings ---
# Use copy/paste-special (formats then values) from a known-good source cell
# so the destination picks up style 13 + t="s" without mutating the style table.
$zeroSrc = $ws.Cells.Item(14,3)   # C14 = style 13, text "0"
$dashSrc = $ws.Cells.Item(14,5)   # E14 = style 13, text "***.*"

$zeroSrc.Copy()
$ws.Cells.Item(15,6).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(15,6).PasteSpecial(-4163)  # xlPasteValues
$zeroSrc.Copy()
$ws.Cells.Item(16,3).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(16,3).PasteSpecial(-4163)  # xlPasteValues
$zeroSrc.Copy()
$ws.Cells.Item(17,3).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(17,3).PasteSpecial(-4163)  # xlPasteValues
$zeroSrc.Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4163)  # xlPasteValues
$zeroSrc.Copy()
$ws.Cells.Item(27,6).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(27,6).PasteSpecial(-4163)  # xlPasteValues
$zeroSrc.Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4122)  # xlPasteFormats
$zeroSrc.Copy()
$ws.Cells.Item(28,4).PasteSpecial(-4163)  # xlPasteValues
$dashSrc.Copy()
$ws.Cells.Item(28,5).PasteSpecial(-4122)  # xlPasteFormats
$dashSrc.Copy()
$ws.Cells.Item(28,5).PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
